$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F12").Value = 1305498328.7699957
$ws.Range("I12").Value = 4188377156

$ws.Range("F13").Value = 325268233.58999997
$ws.Range("I13").Value = 10122006300

$ws.Range("I14").Value = -44319159.289999999

$ws.Range("F16").Value = -53616441.74000001
$ws.Range("I16").Value = -162861893.59999999

$ws.Range("F18").Formula = "=SUM(F12:F17)"

$ws.Range("F19").Value = -412700000
$ws.Range("I19").Value = -1160500000

$ws.Range("F21").Formula = "=SUM(F18:F20)"

$ws.Range("F22").Value = -20015625

$ws.Range("F26").Value = 1026703455.3810816
$ws.Range("I26").Value = 1010658959

$excel.CalculateFullRebuild()
